# step_2 updated and reran
# Append a new slide (DiSCoVER: top drugs, cerebellar stem cell control) at
# the end of the deck. The new slide is a duplicate of the existing
# "DiSCoVER: top drugs (cerebellar stem cell control)" slide (slide 6),
# moved to become the last slide in the deck.

$p = $ppt.ActivePresentation

# Source slide to clone - the existing DiSCoVER top-drugs slide.
$src = $p.Slides.Item(6)

# Duplicate() returns a SlideRange containing the new slide, inserted
# immediately after the source slide.
$dup = $src.Duplicate()
$newSlide = $dup.Item(1)

# Move the freshly duplicated slide to the very end of the deck.
$newSlide.MoveTo($p.Slides.Count)
